# feat: add 2022-Q1 data
#
# The workbook currently has 4 sheets: 2021-Q2, 2021-Q3, 2021-Q4, 总计 (a
# running summary of quarter -> fund-count/holding-value).
#
# This script:
#   1. Turns the existing "总计" sheet into the new "2022-Q1" quarterly
#      holdings sheet (same physical sheet / position as before - so it
#      keeps rId4), and fills it with the 2022-Q1 fund holdings.
#   2. Adds a brand-new "总计" sheet right after it (gets a fresh rId5),
#      and (re)writes the summary table there, with a new 2022-Q1 row
#      prepended above the previously-existing quarters.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: repurpose the current "总计" sheet (4th tab) into "2022-Q1"
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item(4)

# wipe any previous content/formatting on that sheet before reuse
$q1.Cells.Clear()
$q1.Name = "2022-Q1"

# ---------------------------------------------------------------------
# Step 2: insert a fresh "总计" sheet right after it
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q1)
$total.Name = "总计"

# ---------------------------------------------------------------------
# Step 3: populate "2022-Q1" with the quarterly fund-holdings table
# ---------------------------------------------------------------------

# Pull header/index-column/body formatting from an existing quarterly
# sheet (2021-Q3) so the new sheet matches the look of its siblings
# (bold, centered, top-aligned, thin-bordered header + index column).
$fmtSrc = $wb.Worksheets.Item(2)

$fmtSrc.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$fmtSrc.Range("A2:A8").Copy()
$q1.Range("A2:A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# headers
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# D:G hold numeric-looking text (matches the other quarter sheets, where
# these columns are stored as plain text, not numbers) - pre-format the
# block as Text so the values are written verbatim instead of being
# coerced to numbers.
$q1.Range("D2:G5").NumberFormat = "@"

$q1Rows = @(
    @(0, "009098", "景顺长城价值领航两年持有期混合", "11.67", "75.58", "1.92", "0.2241", 8),
    @(1, "005335", "浙商全景消费混合",                 "2.30",  "93.36", "8.23", "0.1893", 3),
    @(2, "010381", "浙商智选价值混合A",                "2.92",  "93.43", "4.81", "0.1405", 8),
    @(3, "010382", "浙商智选价值混合C",                "0.34",  "93.43", "4.81", "0.0164", 8)
)

$r = 2
foreach ($row in $q1Rows) {
    $q1.Range("A$r").Value = $row[0]
    $q1.Range("B$r").Value = $row[1]
    $q1.Range("C$r").Value = $row[2]
    $q1.Range("D$r").Value = $row[3]
    $q1.Range("E$r").Value = $row[4]
    $q1.Range("F$r").Value = $row[5]
    $q1.Range("G$r").Value = $row[6]
    $q1.Range("H$r").Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Step 4: populate the new "总计" sheet with the refreshed summary
# ---------------------------------------------------------------------

$fmtSrc.Range("B1:H1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$fmtSrc.Range("A2:A8").Copy()
$total.Range("A2:A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$totalRows = @(
    @(0, "2022-Q1", 4, 0.57),
    @(1, "2021-Q4", 6, 1.16),
    @(2, "2021-Q3", 7, 1.59),
    @(3, "2021-Q2", 2, 0.14)
)

$r = 2
foreach ($row in $totalRows) {
    $total.Range("A$r").Value = $row[0]
    $total.Range("B$r").Value = $row[1]
    $total.Range("C$r").Value = $row[2]
    $total.Range("D$r").Value = $row[3]
    $r = $r + 1
}

$wb.Worksheets.Item(1).Select()
